# Fruta / hortaliza, semanal
# Insert this week's two new price rows (Kiwi - Vega Monumental Concepción) at the
# top of the data block (row 73), pushing all existing data rows down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current first data row of this block (row 73).
$ws.Rows.Item(73).Insert()
$ws.Rows.Item(73).Insert()

# New row 73: Primera
$ws.Range("A73").Value = 11
$ws.Range("B73").Value = "Vega Monumental Concepción"
$ws.Range("C73").Value = "Bíobío"
$ws.Range("D73").Value = 44775
$ws.Range("E73").Value = 8
$ws.Range("F73").Value = "Fruta"
$ws.Range("G73").Value = 100101
$ws.Range("H73").Value = "Berries"
$ws.Range("I73").Value = 100101007
$ws.Range("J73").Value = "Kiwi"
$ws.Range("K73").Value = "Hayward"
$ws.Range("L73").Value = "Primera"
$ws.Range("M73").Value = 50
$ws.Range("N73").Value = 8000
$ws.Range("O73").Value = 8000
$ws.Range("P73").Value = 8000
$ws.Range("Q73").Value = "$/bandeja 18 kilos"
$ws.Range("R73").Value = "Región de O'Higgins"
$ws.Range("S73").Value = 444
$ws.Range("T73").Value = 18

# New row 74: Segunda
$ws.Range("A74").Value = 11
$ws.Range("B74").Value = "Vega Monumental Concepción"
$ws.Range("C74").Value = "Bíobío"
$ws.Range("D74").Value = 44775
$ws.Range("E74").Value = 8
$ws.Range("F74").Value = "Fruta"
$ws.Range("G74").Value = 100101
$ws.Range("H74").Value = "Berries"
$ws.Range("I74").Value = 100101007
$ws.Range("J74").Value = "Kiwi"
$ws.Range("K74").Value = "Hayward"
$ws.Range("L74").Value = "Segunda"
$ws.Range("M74").Value = 50
$ws.Range("N74").Value = 7000
$ws.Range("O74").Value = 7000
$ws.Range("P74").Value = 7000
$ws.Range("Q74").Value = "$/bandeja 18 kilos"
$ws.Range("R74").Value = "Región de O'Higgins"
$ws.Range("S74").Value = 389
$ws.Range("T74").Value = 18
